$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 49314
$ws.Range("J75").Value = 49314
$ws.Range("L75").Value = 49314
$ws.Range("N75").Value = -51186

$ws.Range("H78").Value = 49314
$ws.Range("J78").Value = 49314
$ws.Range("L78").Value = 147942
$ws.Range("N78").Value = -157302

$ws.Range("H92").Value = 362.14285
$ws.Range("I92").Value = 356.45456
$ws.Range("J92").Value = 383
$ws.Range("K92").Value = 356.45456
$ws.Range("L92").Value = 383
$ws.Range("M92").Value = 891.54544
$ws.Range("N92").Value = -2879

$ws.Range("H112").Value = 1190.2
$ws.Range("J112").Value = 1584.409
$ws.Range("L112").Value = 4753.227000000001
$ws.Range("N112").Value = -6969.227000000001

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0

$ws.Range("H137").Value = 2496.7778
$ws.Range("I137").Value = 1124
$ws.Range("J137").Value = 3595
$ws.Range("K137").Value = 3372
$ws.Range("L137").Value = 10785
$ws.Range("M137").Value = -822
$ws.Range("N137").Value = -15885

$ws.Range("H138").Value = 2655.07
$ws.Range("J138").Value = 2999.0715
$ws.Range("L138").Value = 8997.2145
$ws.Range("N138").Value = -19277.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5465.4
$ws.Range("I32").Value = 5465.4
$ws.Range("K32").Value = 5465.4
$ws.Range("M32").Value = -5178.4

$ws.Range("H45").Value = 1862.25
$ws.Range("I45").Value = 1862.25
$ws.Range("K45").Value = 1862.25
$ws.Range("M45").Value = -1485.25

$ws.Range("H61").Value = 2106.6875
$ws.Range("I61").Value = 1459.421
$ws.Range("J61").Value = 3052.6924
$ws.Range("K61").Value = 1459.421
$ws.Range("L61").Value = 3052.6924
$ws.Range("M61").Value = -1247.421
$ws.Range("N61").Value = -3476.6924

$ws.Range("H74").Value = 22215776
$ws.Range("I74").Value = 24992274
$ws.Range("K74").Value = 24992274
$ws.Range("M74").Value = -24991400

$ws.Range("H77").Value = 22215776
$ws.Range("I77").Value = 24992274
$ws.Range("K77").Value = 124961370
$ws.Range("M77").Value = -124957002

$ws.Range("H132").Value = 2827.1765
$ws.Range("I132").Value = 1758.25
$ws.Range("K132").Value = 5274.75
$ws.Range("M132").Value = -2744.75

$ws.Range("H136").Value = 2106.6875
$ws.Range("I136").Value = 1459.421
$ws.Range("J136").Value = 3052.6924
$ws.Range("K136").Value = 4378.263
$ws.Range("L136").Value = 9158.0772
$ws.Range("M136").Value = -1828.263
$ws.Range("N136").Value = -14258.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 632.36365
$ws.Range("I22").Value = 617.6667
$ws.Range("K22").Value = 617.6667
$ws.Range("M22").Value = -444.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0

$ws.Range("H7").Value = 424.22223
$ws.Range("I7").Value = 269.66666
$ws.Range("K7").Value = 269.66666
$ws.Range("M7").Value = -156.66666

$ws.Range("H16").Value = 1143.091
$ws.Range("I16").Value = 1222.625
$ws.Range("J16").Value = 931
$ws.Range("K16").Value = 1222.625
$ws.Range("L16").Value = 931
$ws.Range("M16").Value = -935.625
$ws.Range("N16").Value = -1505

$ws.Range("H31").Value = 5232.8335
$ws.Range("I31").Value = 1474.75
$ws.Range("J31").Value = 12749
$ws.Range("K31").Value = 1474.75
$ws.Range("L31").Value = 12749
$ws.Range("M31").Value = -1179.75
$ws.Range("N31").Value = -13339

$ws.Range("H34").Value = 5232.8335
$ws.Range("I34").Value = 1474.75
$ws.Range("J34").Value = 12749
$ws.Range("K34").Value = 1474.75
$ws.Range("L34").Value = 12749
$ws.Range("M34").Value = -1272.75
$ws.Range("N34").Value = -13153

$ws.Range("H113").Value = 1143.091
$ws.Range("I113").Value = 1222.625
$ws.Range("J113").Value = 931
$ws.Range("K113").Value = 1222.625
$ws.Range("L113").Value = 931
$ws.Range("M113").Value = 947.375
$ws.Range("N113").Value = -5271

$ws.Range("H115").Value = 49999
$ws.Range("J115").Value = 49999
$ws.Range("L115").Value = 49999
$ws.Range("N115").Value = -52349

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2645680
$ws.Range("J2").Value = 162
$ws.Range("L2").Value = 972
$ws.Range("N2").Value = -1198

$ws.Range("H123").Value = 6000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 6000
$ws.Range("K123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("M123").Value = 18000
$ws.Range("N123").Value = -22900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6153.6665
$ws.Range("I80").Value = 5795.4
$ws.Range("J80").Value = 6601.5
$ws.Range("K80").Value = 5795.4
$ws.Range("L80").Value = 6601.5
$ws.Range("M80").Value = -4797.4
$ws.Range("N80").Value = -8597.5

$ws.Range("H83").Value = 6153.6665
$ws.Range("I83").Value = 5795.4
$ws.Range("J83").Value = 6601.5
$ws.Range("K83").Value = 28977
$ws.Range("L83").Value = 33007.5
$ws.Range("M83").Value = -23985
$ws.Range("N83").Value = -42991.5

$ws.Range("H126").Value = 1303.0769
$ws.Range("I126").Value = 993
$ws.Range("J126").Value = 1799.2
$ws.Range("K126").Value = 2979
$ws.Range("L126").Value = 5397.6
$ws.Range("M126").Value = -509
$ws.Range("N126").Value = -10337.6

$ws.Range("H132").Value = 3021.4666
$ws.Range("I132").Value = 1263
$ws.Range("J132").Value = 3461.0833
$ws.Range("K132").Value = 3789
$ws.Range("L132").Value = 10383.2499
$ws.Range("M132").Value = -1259
$ws.Range("N132").Value = -15443.2499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2224

$ws.Range("H55").Value = 133.66667
$ws.Range("J55").Value = 200
$ws.Range("L55").Value = 200
$ws.Range("N55").Value = -546

$ws.Range("H61").Value = 1427.5714
$ws.Range("I61").Value = 838.8
$ws.Range("K61").Value = 838.8
$ws.Range("M61").Value = -636.8

$ws.Range("H68").Value = 2943.25
$ws.Range("I68").Value = 2936.7144
$ws.Range("K68").Value = 2936.7144
$ws.Range("M68").Value = -2187.7144

$ws.Range("H71").Value = 2943.25
$ws.Range("I71").Value = 2936.7144
$ws.Range("K71").Value = 14683.572
$ws.Range("M71").Value = -10939.572

$ws.Range("H113").Value = 1427.5714
$ws.Range("I113").Value = 838.8
$ws.Range("K113").Value = 838.8
$ws.Range("M113").Value = 1331.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61351
$ws.Range("J46").Value = 61351
$ws.Range("L46").Value = 61351
$ws.Range("N46").Value = -61813

$ws.Range("H62").Value = 102478.2
$ws.Range("J62").Value = 3200
$ws.Range("L62").Value = 3200
$ws.Range("N62").Value = -4448

$ws.Range("H65").Value = 102478.2
$ws.Range("J65").Value = 3200
$ws.Range("L65").Value = 16000
$ws.Range("N65").Value = -22240

$ws.Range("H81").Value = 15123.5
$ws.Range("J81").Value = 19497
$ws.Range("L81").Value = 38994
$ws.Range("N81").Value = -41116

$ws.Range("H84").Value = 15123.5
$ws.Range("J84").Value = 19497
$ws.Range("L84").Value = 194970
$ws.Range("N84").Value = -205578

$ws.Range("H134").Value = 61351
$ws.Range("J134").Value = 61351
$ws.Range("L134").Value = 184053
$ws.Range("N134").Value = -189123
